$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 168.5
$ws.Cells.Item(5, 10).Value = 79
$ws.Cells.Item(5, 12).Value = 79
$ws.Cells.Item(5, 14).Value = -309
$ws.Cells.Item(33, 8).Value = 421.4762
$ws.Cells.Item(33, 9).Value = 248.76923
$ws.Cells.Item(33, 10).Value = 702.125
$ws.Cells.Item(33, 11).Value = 248.76923
$ws.Cells.Item(33, 12).Value = 702.125
$ws.Cells.Item(33, 13).Value = -19.76922999999999
$ws.Cells.Item(33, 14).Value = -1160.125
$ws.Cells.Item(41, 8).Value = 62914.188
$ws.Cells.Item(41, 9).Value = 299.5
$ws.Cells.Item(41, 10).Value = 83785.75
$ws.Cells.Item(41, 11).Value = 299.5
$ws.Cells.Item(41, 12).Value = 83785.75
$ws.Cells.Item(41, 13).Value = 140.5
$ws.Cells.Item(41, 14).Value = -84665.75
$ws.Cells.Item(55, 8).Value = 553.625
$ws.Cells.Item(55, 9).Value = 495.8
$ws.Cells.Item(55, 10).Value = 650
$ws.Cells.Item(55, 11).Value = 495.8
$ws.Cells.Item(55, 12).Value = 650
$ws.Cells.Item(55, 13).Value = -281.8
$ws.Cells.Item(55, 14).Value = -1078
$ws.Cells.Item(100, 8).Value = 6712.1113
$ws.Cells.Item(100, 9).Value = 9284
$ws.Cells.Item(100, 11).Value = 9284
$ws.Cells.Item(100, 13).Value = -8743
$ws.Cells.Item(107, 8).Value = 570.6667
$ws.Cells.Item(107, 9).Value = 570.6667
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 570.6667
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 1349.3333
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 1755.4
$ws.Cells.Item(137, 9).Value = 1355.6875
$ws.Cells.Item(137, 10).Value = 2466
$ws.Cells.Item(137, 11).Value = 4067.0625
$ws.Cells.Item(137, 12).Value = 7398
$ws.Cells.Item(137, 13).Value = -1517.0625
$ws.Cells.Item(137, 14).Value = -12498

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4092.7097
$ws.Cells.Item(32, 9).Value = 4034.0386
$ws.Cells.Item(32, 10).Value = 4397.8
$ws.Cells.Item(32, 11).Value = 4034.0386
$ws.Cells.Item(32, 12).Value = 4397.8
$ws.Cells.Item(32, 13).Value = -3747.0386
$ws.Cells.Item(32, 14).Value = -4971.8
$ws.Cells.Item(61, 8).Value = 18520722
$ws.Cells.Item(61, 9).Value = 23811068
$ws.Cells.Item(61, 11).Value = 23811068
$ws.Cells.Item(61, 13).Value = -23810856
$ws.Cells.Item(74, 8).Value = 1937.1923
$ws.Cells.Item(74, 10).Value = 2951.8333
$ws.Cells.Item(74, 12).Value = 2951.8333
$ws.Cells.Item(74, 14).Value = -4699.8333
$ws.Cells.Item(77, 8).Value = 1937.1923
$ws.Cells.Item(77, 10).Value = 2951.8333
$ws.Cells.Item(77, 12).Value = 14759.1665
$ws.Cells.Item(77, 14).Value = -23495.1665
$ws.Cells.Item(102, 8).Value = 4825915
$ws.Cells.Item(102, 9).Value = 5349552
$ws.Cells.Item(102, 11).Value = 5349552
$ws.Cells.Item(102, 13).Value = -5347930
$ws.Cells.Item(122, 8).Value = 10103982
$ws.Cells.Item(122, 9).Value = 12823251
$ws.Cells.Item(122, 10).Value = 3840.1428
$ws.Cells.Item(122, 11).Value = 38469753
$ws.Cells.Item(122, 12).Value = 11520.4284
$ws.Cells.Item(122, 13).Value = -38467303
$ws.Cells.Item(122, 14).Value = -16420.4284
$ws.Cells.Item(136, 8).Value = 18520722
$ws.Cells.Item(136, 9).Value = 23811068
$ws.Cells.Item(136, 11).Value = 71433204
$ws.Cells.Item(136, 13).Value = -71430654

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 29963.857
$ws.Cells.Item(20, 9).Value = 99999
$ws.Cells.Item(20, 10).Value = 1949.8
$ws.Cells.Item(20, 11).Value = 99999
$ws.Cells.Item(20, 12).Value = 1949.8
$ws.Cells.Item(20, 13).Value = -99752
$ws.Cells.Item(20, 14).Value = -2443.8
$ws.Cells.Item(86, 8).Value = 5863.75
$ws.Cells.Item(86, 9).Value = 6136.6
$ws.Cells.Item(86, 11).Value = 6136.6
$ws.Cells.Item(86, 13).Value = -5013.6
$ws.Cells.Item(89, 8).Value = 5863.75
$ws.Cells.Item(89, 9).Value = 6136.6
$ws.Cells.Item(89, 11).Value = 30683
$ws.Cells.Item(89, 13).Value = -25067
$ws.Cells.Item(94, 8).Value = 1731.8889
$ws.Cells.Item(94, 9).Value = 1324.625
$ws.Cells.Item(94, 10).Value = 4990
$ws.Cells.Item(94, 11).Value = 1324.625
$ws.Cells.Item(94, 12).Value = 4990
$ws.Cells.Item(94, 13).Value = -873.625
$ws.Cells.Item(94, 14).Value = -5892
$ws.Cells.Item(99, 8).Value = 3592
$ws.Cells.Item(99, 9).Value = 3592
$ws.Cells.Item(99, 11).Value = 3592
$ws.Cells.Item(99, 13).Value = -2094
$ws.Cells.Item(105, 8).Value = 1907.8182
$ws.Cells.Item(105, 10).Value = 2070.3333
$ws.Cells.Item(105, 12).Value = 2070.3333
$ws.Cells.Item(105, 14).Value = -5564.3333
$ws.Cells.Item(134, 8).Value = 2012.8136
$ws.Cells.Item(134, 9).Value = 1938.3062
$ws.Cells.Item(134, 10).Value = 2377.9
$ws.Cells.Item(134, 11).Value = 5814.9186
$ws.Cells.Item(134, 12).Value = 7133.700000000001
$ws.Cells.Item(134, 13).Value = -3279.9186
$ws.Cells.Item(134, 14).Value = -12203.7

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(44, 8).Value = 10000
$ws.Cells.Item(44, 9).Value = 10000
$ws.Cells.Item(44, 11).Value = 10000
$ws.Cells.Item(44, 13).Value = -9558
$ws.Cells.Item(47, 8).Value = 42499.75
$ws.Cells.Item(47, 10).Value = 50000
$ws.Cells.Item(47, 12).Value = 50000
$ws.Cells.Item(47, 14).Value = -51132
$ws.Cells.Item(99, 8).Value = 2687.5
$ws.Cells.Item(99, 10).Value = 3000
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 14).Value = -5996
$ws.Cells.Item(105, 8).Value = 1715.5
$ws.Cells.Item(105, 9).Value = 1658.6
$ws.Cells.Item(105, 10).Value = 2000
$ws.Cells.Item(105, 11).Value = 1658.6
$ws.Cells.Item(105, 12).Value = 2000
$ws.Cells.Item(105, 13).Value = 88.40000000000009
$ws.Cells.Item(105, 14).Value = -5494
$ws.Cells.Item(122, 8).Value = 2312.0625
$ws.Cells.Item(122, 9).Value = 2184.4614
$ws.Cells.Item(122, 10).Value = 2865
$ws.Cells.Item(122, 11).Value = 6553.3842
$ws.Cells.Item(122, 12).Value = 8595
$ws.Cells.Item(122, 13).Value = -4103.3842
$ws.Cells.Item(122, 14).Value = -13495
$ws.Cells.Item(126, 8).Value = 2687.5
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 14).Value = -13940
$ws.Cells.Item(132, 8).Value = 3139.8667
$ws.Cells.Item(132, 9).Value = 3199.8635
$ws.Cells.Item(132, 10).Value = 2974.875
$ws.Cells.Item(132, 11).Value = 9599.5905
$ws.Cells.Item(132, 12).Value = 8924.625
$ws.Cells.Item(132, 13).Value = -7069.5905
$ws.Cells.Item(132, 14).Value = -13984.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 493.3793
$ws.Cells.Item(2, 10).Value = 98.75
$ws.Cells.Item(2, 12).Value = 592.5
$ws.Cells.Item(2, 14).Value = -818.5
$ws.Cells.Item(29, 8).Value = 800
$ws.Cells.Item(29, 9).Value = 800
$ws.Cells.Item(29, 11).Value = 2400
$ws.Cells.Item(29, 13).Value = -2123
$ws.Cells.Item(112, 8).Value = 3692
$ws.Cells.Item(112, 9).Value = 3692
$ws.Cells.Item(112, 11).Value = 11076
$ws.Cells.Item(112, 13).Value = -9968
$ws.Cells.Item(140, 8).Value = 1356.6666
$ws.Cells.Item(140, 9).Value = 1356.6666
$ws.Cells.Item(140, 11).Value = 4069.9998
$ws.Cells.Item(140, 13).Value = 1110.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10369.615
$ws.Cells.Item(70, 10).Value = 10429.429
$ws.Cells.Item(70, 12).Value = 10429.429
$ws.Cells.Item(70, 14).Value = -10969.429
$ws.Cells.Item(73, 8).Value = 10369.615
$ws.Cells.Item(73, 10).Value = 10429.429
$ws.Cells.Item(73, 12).Value = 10429.429
$ws.Cells.Item(73, 14).Value = -12301.429
$ws.Cells.Item(126, 8).Value = 14951.375
$ws.Cells.Item(126, 9).Value = 26952.75
$ws.Cells.Item(126, 10).Value = 2950
$ws.Cells.Item(126, 11).Value = 80858.25
$ws.Cells.Item(126, 12).Value = 8850
$ws.Cells.Item(126, 13).Value = -78388.25
$ws.Cells.Item(126, 14).Value = -13790

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 125003064
$ws.Cells.Item(7, 9).Value = 250001120
$ws.Cells.Item(7, 11).Value = 250001120
$ws.Cells.Item(7, 13).Value = -250001008
$ws.Cells.Item(22, 8).Value = 1042.2667
$ws.Cells.Item(22, 9).Value = 1052.5
$ws.Cells.Item(22, 10).Value = 1001.3333
$ws.Cells.Item(22, 11).Value = 1052.5
$ws.Cells.Item(22, 12).Value = 1001.3333
$ws.Cells.Item(22, 13).Value = -757.5
$ws.Cells.Item(22, 14).Value = -1591.3333
$ws.Cells.Item(27, 8).Value = 1042.2667
$ws.Cells.Item(27, 9).Value = 1052.5
$ws.Cells.Item(27, 10).Value = 1001.3333
$ws.Cells.Item(27, 11).Value = 1052.5
$ws.Cells.Item(27, 12).Value = 1001.3333
$ws.Cells.Item(27, 13).Value = -945.5
$ws.Cells.Item(27, 14).Value = -1215.3333
$ws.Cells.Item(46, 8).Value = 2817.5293
$ws.Cells.Item(46, 10).Value = 3722.3333
$ws.Cells.Item(46, 12).Value = 3722.3333
$ws.Cells.Item(46, 14).Value = -4098.3333
$ws.Cells.Item(55, 8).Value = 408.44446
$ws.Cells.Item(55, 9).Value = 458.27274
$ws.Cells.Item(55, 10).Value = 330.14285
$ws.Cells.Item(55, 11).Value = 458.27274
$ws.Cells.Item(55, 12).Value = 330.14285
$ws.Cells.Item(55, 13).Value = -285.27274
$ws.Cells.Item(55, 14).Value = -676.14285
$ws.Cells.Item(93, 8).Value = 1664.8462
$ws.Cells.Item(93, 9).Value = 1624.1818
$ws.Cells.Item(93, 11).Value = 1624.1818
$ws.Cells.Item(93, 13).Value = -376.1818000000001
$ws.Cells.Item(126, 8).Value = 125003064
$ws.Cells.Item(126, 9).Value = 250001120
$ws.Cells.Item(126, 11).Value = 750003360
$ws.Cells.Item(126, 13).Value = -750000890

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1849.6666
$ws.Cells.Item(126, 9).Value = 1774.5
$ws.Cells.Item(126, 11).Value = 5323.5
$ws.Cells.Item(126, 13).Value = -2853.5
$ws.Cells.Item(132, 8).Value = 5055.316
$ws.Cells.Item(132, 9).Value = 5658
$ws.Cells.Item(132, 10).Value = 3749.5
$ws.Cells.Item(132, 11).Value = 16974
$ws.Cells.Item(132, 12).Value = 11248.5
$ws.Cells.Item(132, 13).Value = -14444
$ws.Cells.Item(132, 14).Value = -16308.5
$ws.Cells.Item(136, 8).Value = 3198.6428
$ws.Cells.Item(136, 9).Value = 1670
$ws.Cells.Item(136, 11).Value = 5010
$ws.Cells.Item(136, 13).Value = -2460
